# Update "想去人数" (F column) values in 展览 and 全部类型 sheets
# to reflect new generated output data.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates (row -> new F value)
$wsExhibit.Range("F2").Value = 601
$wsExhibit.Range("F3").Value = 494
$wsExhibit.Range("F4").Value = 1276
$wsExhibit.Range("F6").Value = 14201
$wsExhibit.Range("F7").Value = 15983
$wsExhibit.Range("F24").Value = 6357
$wsExhibit.Range("F25").Value = 966
$wsExhibit.Range("F26").Value = 1106
$wsExhibit.Range("F27").Value = 5636
$wsExhibit.Range("F28").Value = 84
$wsExhibit.Range("F30").Value = 144
$wsExhibit.Range("F31").Value = 4636

# 全部类型 sheet updates (row -> new F value)
$wsAll.Range("F2").Value = 601
$wsAll.Range("F3").Value = 494
$wsAll.Range("F4").Value = 1276
$wsAll.Range("F6").Value = 14201
$wsAll.Range("F7").Value = 15983
$wsAll.Range("F25").Value = 6357
$wsAll.Range("F26").Value = 966
$wsAll.Range("F27").Value = 1106
$wsAll.Range("F29").Value = 5636
$wsAll.Range("F30").Value = 84
$wsAll.Range("F32").Value = 144
$wsAll.Range("F33").Value = 4636
